$d = $word.ActiveDocument

# --- 1) Shrink ", portador do CNPJ:" down to ", " and insert "inscrito no CNPJ:" after it ---
$rLead = $d.Content
$null = $rLead.Find.Execute(", portador do CNPJ:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rLead.Text = ", "
$rLead.Collapse(0)
$rLead.InsertAfter("inscrito no CNPJ:")

# --- 2) Replace the bold " #CNPJ" run with the longer sentence, keeping only "#CNPJ" bold ---
#     Empty the existing bold run first so the new text we type inherits the surrounding
#     (non-bold) formatting instead of carrying the old run's bold flag along with it.
$rCnpj = $d.Content
$null = $rCnpj.Find.Execute(" #CNPJ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rCnpj.Text = ""
$insPoint = $rCnpj.Start
$rCnpj.InsertAfter(" #CNPJ e com sede na #END_EMPRESA cep: #CP_EMPRESA")

# Re-bold just the "#CNPJ" placeholder portion (" #CNPJ" = 6 characters).
$rBold = $d.Range($insPoint, $insPoint + 6)
$rBold.Bold = 1
